# Commit: "Generate Report for Archive"
#
# The underlying OOXML diff for this commit only touches xl/sharedStrings.xml:
# a new, distinct string value ("In Translation") is appended to the shared
# string table (count/uniqueCount 56 -> 57), and every worksheet <c t="s"><v>
# index that pointed at or past the insertion point is bumped by exactly one
# slot to keep pointing at the same text as before.
#
# Concretely: every single cell on every sheet (Overview, zh-cn, de-de)
# resolves, before and after the commit, to the exact same displayed text --
# "Ready for handoff" is still "Ready for handoff", every date/guid/path
# string is unchanged, etc. The new shared-string entry is never referenced
# by any cell in this snapshot; it is simply pre-registered in the status
# vocabulary used by the report generator (presumably so a future run that
# *does* hit that status doesn't need to grow the table again), which is
# exactly the kind of bookkeeping a localization-status report generator
# tool does on a refresh/regeneration pass ("Generate Report for Archive")
# without any of the actual report data changing.
#
# Shared-string allocation in the workbook's string pool is an internal
# serialization detail of the writer, not something exposed on the Excel
# object model (there's no Application/Workbook method to poke the shared
# string table directly, and real Excel manages it transparently too) --
# so the faithful, object-model-level replay of this commit is to leave
# every sheet's cell values exactly as they already are: re-generating the
# report against unchanged source data reproduces the same values.

$wb = $excel.ActiveWorkbook

# Touch each sheet (mirrors the report generator re-visiting every sheet
# when it regenerates/archives the workbook) without altering any cell's
# value, formula, or formatting.
foreach ($ws in $wb.Worksheets) {
    $ws.Calculate()
}

$wb.Save()
